$wb = $excel.ActiveWorkbook

$values = @{
    "N-Dense" = "5.48"
    "N-Type" = "5.89"
    "N-type Wafer" = "1.19"
    "Cell Topcon 183mm" = "0.29"
    "Module Topcon 183mm" = "0.1"
    "Silver Rear_side" = "5,195"
    "Silver Busbar front-side" = "7,777"
    "Silver finger front-side" = "7,827"
    "USD_CNY" = "7.3048"
}

foreach ($sheetName in $values.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A4:B4").NumberFormat = "@"
    $ws.Cells.Item(4, 1).Value = "2025-03-05"
    $ws.Cells.Item(4, 2).Value = $values[$sheetName]
}
